# ESCALETA_ MA_09_13_CO.xlsx - "Ajuste de escaleta mat 9 tema 13"
#
# The sheet "Hoja2" (the first / active sheet) has an AutoFilter on column P
# that hides rows whose P value isn't "SI". A handful of rows that were
# hidden by that filter (3, 10, 24 and the whole trailing block 254-279),
# plus two already-visible rows (27, 28), had their row height manually
# adjusted (auto-fit) and - for the previously filtered-out rows - were
# unhidden again. The active selection also moved from U23 to J256.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Row -> new (auto-fit) row height, taken from the target workbook.
$rowHeights = @{
    3   = 16.5
    10  = 14.25
    24  = 14.25
    27  = 15.75
    28  = 28.5
    254 = 20.25
    255 = 20.25
    256 = 30.75
    257 = 29.25
    258 = 20.25
    259 = 18
    260 = 27.75
    261 = 21
    262 = 27.75
    263 = 23.25
    264 = 25.5
    265 = 21
    266 = 26.25
    267 = 19.5
    268 = 21.75
    269 = 27
    270 = 26.25
    271 = 21.75
    272 = 18
    273 = 18
    274 = 24
    275 = 13.5
    276 = 29.25
    277 = 22.5
    278 = 48
    279 = 26.25
}

foreach ($r in $rowHeights.Keys) {
    $row = $ws.Rows.Item($r)
    $row.Hidden = $false
    $row.RowHeight = $rowHeights[$r]
}

# Move the active selection to where the author left it.
$ws.Range("J256").Select()
